$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G = "Recorded By" contains comma-separated lists of recorder
# names/emails. This edit re-orders the entries in specific, known
# combinations (same set of values, different order) to match the
# upstream sync of attendance_reports. String comparisons must be
# case-sensitive (Ordinal) because "system" and "System" both appear
# as distinct tokens.

$ordinal = [System.StringComparison]::Ordinal

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    $parts = $val -split ', '

    $newVal = $val

    if ($parts.Length -eq 2 -and $parts[0].Equals('dnasr281@gmail.com', $ordinal) -and $parts[1].Equals('System', $ordinal)) {
        $newVal = 'System, dnasr281@gmail.com'
    }
    elseif ($parts.Length -eq 2 -and $parts[0].Equals('dnasr281@gmail.com', $ordinal) -and $parts[1].Equals('admin@admin.com', $ordinal)) {
        $newVal = 'admin@admin.com, dnasr281@gmail.com'
    }
    elseif ($parts.Length -eq 3 -and $parts[0].Equals('backup@backdoor.com', $ordinal) -and $parts[1].Equals('system', $ordinal) -and $parts[2].Equals('System', $ordinal)) {
        $newVal = 'backup@backdoor.com, System, system'
    }

    if (-not $newVal.Equals($val, $ordinal)) {
        $cell.Value = $newVal
    }
}
